$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 859.2
$ws.Cells.Item(19, 9).Value = 536.75
$ws.Cells.Item(19, 10).Value = 1074.1666
$ws.Cells.Item(19, 11).Value = 536.75
$ws.Cells.Item(19, 12).Value = 1074.1666
$ws.Cells.Item(19, 13).Value = -361.75
$ws.Cells.Item(19, 14).Value = -1424.1666
$ws.Cells.Item(137, 8).Value = 1860.98
$ws.Cells.Item(137, 9).Value = 1376.5588
$ws.Cells.Item(137, 10).Value = 2890.375
$ws.Cells.Item(137, 11).Value = 4129.6764
$ws.Cells.Item(137, 12).Value = 8671.125
$ws.Cells.Item(137, 13).Value = -1579.6764
$ws.Cells.Item(137, 14).Value = -13771.125
$ws.Cells.Item(141, 8).Value = 2908.8774
$ws.Cells.Item(141, 9).Value = 1674.25
$ws.Cells.Item(141, 10).Value = 5232.8823
$ws.Cells.Item(141, 11).Value = 5022.75
$ws.Cells.Item(141, 12).Value = 15698.6469
$ws.Cells.Item(141, 13).Value = 157.25
$ws.Cells.Item(141, 14).Value = -26058.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17701.568
$ws.Cells.Item(32, 9).Value = 17789.598
$ws.Cells.Item(32, 10).Value = 16859
$ws.Cells.Item(32, 11).Value = 17789.598
$ws.Cells.Item(32, 12).Value = 16859
$ws.Cells.Item(32, 13).Value = -17502.598
$ws.Cells.Item(32, 14).Value = -17433
$ws.Cells.Item(45, 8).Value = 1335.678
$ws.Cells.Item(45, 9).Value = 1266.3959
$ws.Cells.Item(45, 10).Value = 1638
$ws.Cells.Item(45, 11).Value = 1266.3959
$ws.Cells.Item(45, 12).Value = 1638
$ws.Cells.Item(45, 13).Value = -889.3959
$ws.Cells.Item(45, 14).Value = -2392
$ws.Cells.Item(61, 8).Value = 5101.1934
$ws.Cells.Item(61, 9).Value = 3075.7693
$ws.Cells.Item(61, 10).Value = 15633.4
$ws.Cells.Item(61, 11).Value = 3075.7693
$ws.Cells.Item(61, 12).Value = 15633.4
$ws.Cells.Item(61, 13).Value = -2863.7693
$ws.Cells.Item(61, 14).Value = -16057.4
$ws.Cells.Item(74, 8).Value = 3565.5
$ws.Cells.Item(74, 9).Value = 1491.7778
$ws.Cells.Item(74, 11).Value = 1491.7778
$ws.Cells.Item(74, 13).Value = -617.7778000000001
$ws.Cells.Item(77, 8).Value = 3565.5
$ws.Cells.Item(77, 9).Value = 1491.7778
$ws.Cells.Item(77, 11).Value = 7458.889
$ws.Cells.Item(77, 13).Value = -3090.889
$ws.Cells.Item(110, 8).Value = 1851.4445
$ws.Cells.Item(110, 9).Value = 1733.25
$ws.Cells.Item(110, 10).Value = 1946
$ws.Cells.Item(110, 11).Value = 1733.25
$ws.Cells.Item(110, 12).Value = 1946
$ws.Cells.Item(110, 13).Value = 311.75
$ws.Cells.Item(110, 14).Value = -6036
$ws.Cells.Item(122, 8).Value = 2149.4285
$ws.Cells.Item(122, 9).Value = 1687.7142
$ws.Cells.Item(122, 10).Value = 2611.1428
$ws.Cells.Item(122, 11).Value = 5063.142599999999
$ws.Cells.Item(122, 12).Value = 7833.428400000001
$ws.Cells.Item(122, 13).Value = -2613.142599999999
$ws.Cells.Item(122, 14).Value = -12733.4284
$ws.Cells.Item(124, 8).Value = 15071.5
$ws.Cells.Item(124, 10).Value = 15071.5
$ws.Cells.Item(124, 12).Value = 15071.5
$ws.Cells.Item(124, 14).Value = -24891.5
$ws.Cells.Item(125, 8).Value = 54073.57
$ws.Cells.Item(125, 10).Value = 54073.57
$ws.Cells.Item(125, 12).Value = 54073.57
$ws.Cells.Item(125, 14).Value = -63913.57
$ws.Cells.Item(132, 8).Value = 5016.6387
$ws.Cells.Item(132, 9).Value = 1485.05
$ws.Cells.Item(132, 10).Value = 9431.125
$ws.Cells.Item(132, 11).Value = 4455.15
$ws.Cells.Item(132, 12).Value = 28293.375
$ws.Cells.Item(132, 13).Value = -1925.15
$ws.Cells.Item(132, 14).Value = -33353.375
$ws.Cells.Item(136, 8).Value = 5101.1934
$ws.Cells.Item(136, 9).Value = 3075.7693
$ws.Cells.Item(136, 10).Value = 15633.4
$ws.Cells.Item(136, 11).Value = 9227.3079
$ws.Cells.Item(136, 12).Value = 46900.2
$ws.Cells.Item(136, 13).Value = -6677.3079
$ws.Cells.Item(136, 14).Value = -52000.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(54, 8).Value = 15024.454
$ws.Cells.Item(54, 9).Value = 2394.3333
$ws.Cells.Item(54, 10).Value = 19760.75
$ws.Cells.Item(54, 11).Value = 2394.3333
$ws.Cells.Item(54, 12).Value = 19760.75
$ws.Cells.Item(54, 13).Value = -1910.3333
$ws.Cells.Item(54, 14).Value = -20728.75
$ws.Cells.Item(99, 8).Value = 1528
$ws.Cells.Item(99, 9).Value = 1309.909
$ws.Cells.Item(99, 11).Value = 1309.909
$ws.Cells.Item(99, 13).Value = 188.0909999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1834.4333
$ws.Cells.Item(31, 9).Value = 1261.5686
$ws.Cells.Item(31, 10).Value = 5080.6665
$ws.Cells.Item(31, 11).Value = 1261.5686
$ws.Cells.Item(31, 12).Value = 5080.6665
$ws.Cells.Item(31, 13).Value = -966.5686000000001
$ws.Cells.Item(31, 14).Value = -5670.6665
$ws.Cells.Item(34, 8).Value = 1834.4333
$ws.Cells.Item(34, 9).Value = 1261.5686
$ws.Cells.Item(34, 10).Value = 5080.6665
$ws.Cells.Item(34, 11).Value = 1261.5686
$ws.Cells.Item(34, 12).Value = 5080.6665
$ws.Cells.Item(34, 13).Value = -1059.5686
$ws.Cells.Item(34, 14).Value = -5484.6665
$ws.Cells.Item(132, 8).Value = 2991.2207
$ws.Cells.Item(132, 9).Value = 3251.0566
$ws.Cells.Item(132, 10).Value = 2417.4167
$ws.Cells.Item(132, 11).Value = 9753.1698
$ws.Cells.Item(132, 12).Value = 7252.250100000001
$ws.Cells.Item(132, 13).Value = -7223.1698
$ws.Cells.Item(132, 14).Value = -12312.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1577.7894
$ws.Cells.Item(131, 9).Value = 3972.5
$ws.Cells.Item(131, 10).Value = 1296.0588
$ws.Cells.Item(131, 11).Value = 11917.5
$ws.Cells.Item(131, 12).Value = 3888.1764
$ws.Cells.Item(131, 13).Value = -6877.5
$ws.Cells.Item(131, 14).Value = -13968.1764

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3473.9375
$ws.Cells.Item(102, 9).Value = 2897.3044
$ws.Cells.Item(102, 10).Value = 4947.5557
$ws.Cells.Item(102, 11).Value = 2897.3044
$ws.Cells.Item(102, 12).Value = 4947.5557
$ws.Cells.Item(102, 13).Value = -1275.3044
$ws.Cells.Item(102, 14).Value = -8191.5557
$ws.Cells.Item(113, 8).Value = 3714.2856
$ws.Cells.Item(113, 10).Value = 2760
$ws.Cells.Item(113, 12).Value = 2760
$ws.Cells.Item(113, 14).Value = -7100
$ws.Cells.Item(123, 8).Value = 10326
$ws.Cells.Item(123, 10).Value = 10326
$ws.Cells.Item(123, 12).Value = 10326
$ws.Cells.Item(123, 14).Value = -15226
$ws.Cells.Item(135, 8).Value = 53796.285
$ws.Cells.Item(135, 10).Value = 53796.285
$ws.Cells.Item(135, 12).Value = 53796.285
$ws.Cells.Item(135, 14).Value = -63936.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 25003916
$ws.Cells.Item(14, 10).Value = 5333
$ws.Cells.Item(14, 12).Value = 5333
$ws.Cells.Item(14, 14).Value = -5669
$ws.Cells.Item(126, 8).Value = 1501.8148
$ws.Cells.Item(126, 9).Value = 1412.7
$ws.Cells.Item(126, 10).Value = 1756.4286
$ws.Cells.Item(126, 11).Value = 4238.1
$ws.Cells.Item(126, 12).Value = 5269.2858
$ws.Cells.Item(126, 13).Value = -1768.1
$ws.Cells.Item(126, 14).Value = -10209.2858
$ws.Cells.Item(132, 8).Value = 1152.6538
$ws.Cells.Item(132, 9).Value = 436.17648
$ws.Cells.Item(132, 10).Value = 2506
$ws.Cells.Item(132, 11).Value = 1308.52944
$ws.Cells.Item(132, 12).Value = 7518
$ws.Cells.Item(132, 13).Value = 1221.47056
$ws.Cells.Item(132, 14).Value = -12578
